$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.348.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "1.794.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'225.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").Value = "'0.595"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.11%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'36.16"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.31%  "
$ws.Range("E9").Value = "  -3.89%  "
$ws.Range("D10").Value = "'0.0676"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.49%  "
$ws.Range("D11").Value = "'0.0961"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").Value = "2.053.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("D13").Value = "'11.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.01%  "
$ws.Range("D14").Value = "1.804.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("E15").Value = "  -2.73%  "
$ws.Range("D16").Value = "34.304.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").Value = "'4.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "'68.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "'240.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.40%  "
$ws.Range("D20").Value = "0.0₃0770"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.90%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "'11.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.41%  "
$ws.Range("D23").Value = "'4.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("E24").Value = "  +4.30%  "
$ws.Range("D25").Value = "'170.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.98%  "
$ws.Range("D26").Value = "'7.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.02%  "
$ws.Range("D27").Value = "'17.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("D28").Value = "'0.120"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("D31").Value = "'3.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("D32").Value = "'3.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("D33").Value = "'0.0511"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.56%  "
$ws.Range("E34").Value = "  -4.89%  "
$ws.Range("D35").Value = "1.356.87"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("D36").Value = "'0.645"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.03%  "
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("E38").Value = "  -9.62%  "
$ws.Range("D39").Value = "'0.0184"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.68%  "
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("E41").Value = "  -3.61%  "
$ws.Range("D42").Value = "'80.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("D43").Value = "'0.930"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.31%  "
$ws.Range("E44").Value = "  +5.70%  "
$ws.Range("D45").Value = "'13.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.20%  "
$ws.Range("D46").Value = "'0.0494"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.13%  "
$ws.Range("D47").Value = "1.955.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("D48").Value = "'5.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.49%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "'101.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.42%  "
$ws.Range("D51").Value = "0.0₆0118"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.76%  "
